# Generate Report for Handback
# - Overview sheet: status for "32566931-..." row switches from
#   "Ready for handoff" to "Handback transform failed" (zh-cn + de-de cols)
# - zh-cn / de-de sheets: fill in the "Error Detail" column (P) for that
#   same row with the handback-mismatch message, and widen column P to fit.

$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"

# The "Ready for handoff" status for the 32566931-... file is shared text
# across the Overview sheet's per-language columns AND each language
# sheet's own Status column - update every occurrence so they stay in sync.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# Excel's ColumnWidth (character units) round-trips through its own pixel
# grid and gets ~5/6 of a character added back on save, so asking for the
# saved/stored width of 40 directly overshoots to ~40.83. Back the padding
# out of the request so the persisted <col> width lands on exactly 40.
$targetColWidth = 40 - (5 / 6)

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $newStatus
$zhcn.Columns.Item(16).ColumnWidth = $targetColWidth
$zhcn.Range("P3").Value = "Handback file name: uhgzvj4o.eir is different with handoff file name: 32566931-2754-4927-a9d4-6e9d4ca69873.d509c8dc932db8a76aab0e785c89e7a320b3dcbf.zh-cn."

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $newStatus
$dede.Columns.Item(16).ColumnWidth = $targetColWidth
$dede.Range("P3").Value = "Handback file name: uhgzvj4o.eir is different with handoff file name: 32566931-2754-4927-a9d4-6e9d4ca69873.d509c8dc932db8a76aab0e785c89e7a320b3dcbf.de-de."
